$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D415").Value = "test"
